# "ignore un nexessary files and pull new changes"
#
# Rename the original sheet, add two more worksheets with their data, and
# add a new "runMode" column to the original sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: AddCustomerTest -> addCustomerTest -----------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "addCustomerTest"

# --- Sheet 2: openAccountTest (new) -------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "openAccountTest"

$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A1:B1").Font.Bold = $true

$ws2.Range("A2").Value = "Sohaib Majeed"
$ws2.Range("B2").Value = "Rupee"

$ws2.Columns.Item(1).ColumnWidth = 14.28
$ws2.Columns.Item(2).ColumnWidth = 8.57
$ws2.Range("A2").Select() | Out-Null

# --- Sheet 3: test_suite (new) ------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "test_suite"

$ws3.Range("A1").Value = "TCID"
$ws3.Range("B1").Value = "Runmode"
$ws3.Range("A1:B1").Font.Bold = $true

$ws3.Range("A2").Value = "BankManagerLoginTest"
$ws3.Range("B2").Value = "Y"
$ws3.Range("A3").Value = "AddCustomerTest"
$ws3.Range("B3").Value = "Y"
$ws3.Range("A4").Value = "OpenAccountTest"
$ws3.Range("B4").Value = "Y"

$ws3.Columns.Item(1).ColumnWidth = 21.85
$ws3.Range("B4").Select() | Out-Null

# --- Back to sheet 1: add the new customer rows + runMode column -------
$ws1.Range("A3").Value = "Daniyal"
$ws1.Range("B3").Value = "Ahmed"
$ws1.Range("C3").Value = "4567xy"
$ws1.Range("D3").Value = "Customer added successfully"

$ws1.Range("A4").Value = "Kashan"
$ws1.Range("B4").Value = "Ali"
$ws1.Range("C4").Value = "76yrt"
$ws1.Range("D4").Value = "Customer added successfully"

$ws1.Range("A5").Value = "Usman"
$ws1.Range("B5").Value = "Shabeer"
$ws1.Range("C5").Value = "89rt"
$ws1.Range("D5").Value = "Customer added successfully"

$ws1.Range("E1").Value = "runMode"
$ws1.Range("E1").Font.Bold = $true
$ws1.Range("E2").Value = "Y"
$ws1.Range("E3").Value = "Y"
$ws1.Range("E4").Value = "Y"
$ws1.Range("E5").Value = "Y"

$ws1.Range("E3").Select() | Out-Null

# test_suite ends up the active tab, matching tabSelected/activeTab="2"
$ws3.Activate() | Out-Null
